$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Static values common to every data row in this block (cols A,B,C,E,F,G,I,R)
$colA = 7
$colB = 'Terminal Hortofrutícola Agro Chillán'
$colC = 'Ñuble'
$colE = 16
$colF = 100112021
$colG = 'Ají'
$colI = 'Primera'
$colR = 'Hortaliza'

# Extend dimension by adding three new rows (167-169) with the common static values,
# then overwrite the variable data (D,H,J,K,L,M,N,O,P,Q) for rows 113-169 to match the
# updated weekly price list (row 113 is the new entry, rest shift down by 3).
for ($r = 167; $r -le 169; $r++) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 18).Value = $colR
}

# Row 113
$ws.Cells.Item(113, 4).Value = 45009
$ws.Cells.Item(113, 8).Value = 'Cacho cabra rojo'
$ws.Cells.Item(113, 10).Value = 20
$ws.Cells.Item(113, 11).Value = 15000
$ws.Cells.Item(113, 12).Value = 15000
$ws.Cells.Item(113, 13).Value = 15000
$ws.Cells.Item(113, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(113, 15).Value = 'Región del Maule'
$ws.Cells.Item(113, 16).Value = 600
$ws.Cells.Item(113, 17).Value = 25

# Row 114
$ws.Cells.Item(114, 4).Value = 45009
$ws.Cells.Item(114, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(114, 10).Value = 20
$ws.Cells.Item(114, 11).Value = 15000
$ws.Cells.Item(114, 12).Value = 15000
$ws.Cells.Item(114, 13).Value = 15000
$ws.Cells.Item(114, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(114, 15).Value = 'Región del Maule'
$ws.Cells.Item(114, 16).Value = 600
$ws.Cells.Item(114, 17).Value = 25

# Row 115
$ws.Cells.Item(115, 4).Value = 45009
$ws.Cells.Item(115, 8).Value = 'Cristal'
$ws.Cells.Item(115, 10).Value = 20
$ws.Cells.Item(115, 11).Value = 15000
$ws.Cells.Item(115, 12).Value = 15000
$ws.Cells.Item(115, 13).Value = 15000
$ws.Cells.Item(115, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(115, 15).Value = 'Región del Maule'
$ws.Cells.Item(115, 16).Value = 600
$ws.Cells.Item(115, 17).Value = 25

# Row 116
$ws.Cells.Item(116, 4).Value = 44942
$ws.Cells.Item(116, 8).Value = 'Americana (o)'
$ws.Cells.Item(116, 10).Value = 60
$ws.Cells.Item(116, 11).Value = 13000
$ws.Cells.Item(116, 12).Value = 13500
$ws.Cells.Item(116, 13).Value = 13250
$ws.Cells.Item(116, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(116, 15).Value = 'Región del Maule'
$ws.Cells.Item(116, 16).Value = 883
$ws.Cells.Item(116, 17).Value = 15

# Row 117
$ws.Cells.Item(117, 4).Value = 44232
$ws.Cells.Item(117, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(117, 10).Value = 30
$ws.Cells.Item(117, 11).Value = 12000
$ws.Cells.Item(117, 12).Value = 13000
$ws.Cells.Item(117, 13).Value = 12500
$ws.Cells.Item(117, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(117, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(117, 16).Value = 500
$ws.Cells.Item(117, 17).Value = 25

# Row 118
$ws.Cells.Item(118, 4).Value = 44294
$ws.Cells.Item(118, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(118, 10).Value = 30
$ws.Cells.Item(118, 11).Value = 11500
$ws.Cells.Item(118, 12).Value = 12000
$ws.Cells.Item(118, 13).Value = 11750
$ws.Cells.Item(118, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(118, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(118, 16).Value = 470
$ws.Cells.Item(118, 17).Value = 25

# Row 119
$ws.Cells.Item(119, 4).Value = 44551
$ws.Cells.Item(119, 8).Value = 'Americana (o)'
$ws.Cells.Item(119, 10).Value = 60
$ws.Cells.Item(119, 11).Value = 17500
$ws.Cells.Item(119, 12).Value = 18000
$ws.Cells.Item(119, 13).Value = 17750
$ws.Cells.Item(119, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(119, 15).Value = 'Región del Maule'
$ws.Cells.Item(119, 16).Value = 1183
$ws.Cells.Item(119, 17).Value = 15

# Row 120
$ws.Cells.Item(120, 4).Value = 44798
$ws.Cells.Item(120, 8).Value = 'Americana (o)'
$ws.Cells.Item(120, 10).Value = 60
$ws.Cells.Item(120, 11).Value = 29000
$ws.Cells.Item(120, 12).Value = 30000
$ws.Cells.Item(120, 13).Value = 29500
$ws.Cells.Item(120, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(120, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(120, 16).Value = 1967
$ws.Cells.Item(120, 17).Value = 15

# Row 121
$ws.Cells.Item(121, 4).Value = 44691
$ws.Cells.Item(121, 8).Value = 'Cristal'
$ws.Cells.Item(121, 10).Value = 60
$ws.Cells.Item(121, 11).Value = 24000
$ws.Cells.Item(121, 12).Value = 25000
$ws.Cells.Item(121, 13).Value = 24500
$ws.Cells.Item(121, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(121, 15).Value = 'Región del Maule'
$ws.Cells.Item(121, 16).Value = 980
$ws.Cells.Item(121, 17).Value = 25

# Row 122
$ws.Cells.Item(122, 4).Value = 44637
$ws.Cells.Item(122, 8).Value = 'Americana (o)'
$ws.Cells.Item(122, 10).Value = 80
$ws.Cells.Item(122, 11).Value = 8500
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = 8750
$ws.Cells.Item(122, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(122, 15).Value = 'Región del Maule'
$ws.Cells.Item(122, 16).Value = 583
$ws.Cells.Item(122, 17).Value = 15

# Row 123
$ws.Cells.Item(123, 4).Value = 44222
$ws.Cells.Item(123, 8).Value = 'Americana (o)'
$ws.Cells.Item(123, 10).Value = 75
$ws.Cells.Item(123, 11).Value = 900
$ws.Cells.Item(123, 12).Value = 1000
$ws.Cells.Item(123, 13).Value = 933
$ws.Cells.Item(123, 14).Value = '$/kilo'
$ws.Cells.Item(123, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(123, 16).Value = 933
$ws.Cells.Item(123, 17).Value = 1

# Row 124
$ws.Cells.Item(124, 4).Value = 44244
$ws.Cells.Item(124, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(124, 10).Value = 20
$ws.Cells.Item(124, 11).Value = 13000
$ws.Cells.Item(124, 12).Value = 14000
$ws.Cells.Item(124, 13).Value = 13500
$ws.Cells.Item(124, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(124, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(124, 16).Value = 540
$ws.Cells.Item(124, 17).Value = 25

# Row 125
$ws.Cells.Item(125, 4).Value = 44258
$ws.Cells.Item(125, 8).Value = 'Americana (o)'
$ws.Cells.Item(125, 10).Value = 75
$ws.Cells.Item(125, 11).Value = 600
$ws.Cells.Item(125, 12).Value = 700
$ws.Cells.Item(125, 13).Value = 633
$ws.Cells.Item(125, 14).Value = '$/kilo'
$ws.Cells.Item(125, 15).Value = 'Región del Maule'
$ws.Cells.Item(125, 16).Value = 633
$ws.Cells.Item(125, 17).Value = 1

# Row 126
$ws.Cells.Item(126, 4).Value = 44566
$ws.Cells.Item(126, 8).Value = 'Americana (o)'
$ws.Cells.Item(126, 10).Value = 60
$ws.Cells.Item(126, 11).Value = 20000
$ws.Cells.Item(126, 12).Value = 21000
$ws.Cells.Item(126, 13).Value = 20500
$ws.Cells.Item(126, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(126, 15).Value = 'Región del Maule'
$ws.Cells.Item(126, 16).Value = 1367
$ws.Cells.Item(126, 17).Value = 15

# Row 127
$ws.Cells.Item(127, 4).Value = 44658
$ws.Cells.Item(127, 8).Value = 'Americana (o)'
$ws.Cells.Item(127, 10).Value = 80
$ws.Cells.Item(127, 11).Value = 8500
$ws.Cells.Item(127, 12).Value = 9000
$ws.Cells.Item(127, 13).Value = 8750
$ws.Cells.Item(127, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(127, 15).Value = 'Región del Maule'
$ws.Cells.Item(127, 16).Value = 583
$ws.Cells.Item(127, 17).Value = 15

# Row 128
$ws.Cells.Item(128, 4).Value = 44782
$ws.Cells.Item(128, 8).Value = 'Inferno'
$ws.Cells.Item(128, 10).Value = 60
$ws.Cells.Item(128, 11).Value = 15000
$ws.Cells.Item(128, 12).Value = 16000
$ws.Cells.Item(128, 13).Value = 15500
$ws.Cells.Item(128, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(128, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(128, 16).Value = 1033
$ws.Cells.Item(128, 17).Value = 15

# Row 129
$ws.Cells.Item(129, 4).Value = 44832
$ws.Cells.Item(129, 8).Value = 'Inferno'
$ws.Cells.Item(129, 10).Value = 60
$ws.Cells.Item(129, 11).Value = 24000
$ws.Cells.Item(129, 12).Value = 25000
$ws.Cells.Item(129, 13).Value = 24500
$ws.Cells.Item(129, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(129, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(129, 16).Value = 2450
$ws.Cells.Item(129, 17).Value = 10

# Row 130
$ws.Cells.Item(130, 4).Value = 44586
$ws.Cells.Item(130, 8).Value = 'Americana (o)'
$ws.Cells.Item(130, 10).Value = 60
$ws.Cells.Item(130, 11).Value = 13000
$ws.Cells.Item(130, 12).Value = 14000
$ws.Cells.Item(130, 13).Value = 13500
$ws.Cells.Item(130, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(130, 15).Value = 'Región del Maule'
$ws.Cells.Item(130, 16).Value = 900
$ws.Cells.Item(130, 17).Value = 15

# Row 131
$ws.Cells.Item(131, 4).Value = 44907
$ws.Cells.Item(131, 8).Value = 'Americana (o)'
$ws.Cells.Item(131, 10).Value = 100
$ws.Cells.Item(131, 11).Value = 15500
$ws.Cells.Item(131, 12).Value = 16000
$ws.Cells.Item(131, 13).Value = 15750
$ws.Cells.Item(131, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(131, 15).Value = 'Región del Maule'
$ws.Cells.Item(131, 16).Value = 1050
$ws.Cells.Item(131, 17).Value = 15

# Row 132
$ws.Cells.Item(132, 4).Value = 44979
$ws.Cells.Item(132, 8).Value = 'Americana (o)'
$ws.Cells.Item(132, 10).Value = 60
$ws.Cells.Item(132, 11).Value = 16000
$ws.Cells.Item(132, 12).Value = 17000
$ws.Cells.Item(132, 13).Value = 16500
$ws.Cells.Item(132, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(132, 15).Value = 'Región del Maule'
$ws.Cells.Item(132, 16).Value = 660
$ws.Cells.Item(132, 17).Value = 25

# Row 133
$ws.Cells.Item(133, 4).Value = 44979
$ws.Cells.Item(133, 8).Value = 'Cristal'
$ws.Cells.Item(133, 10).Value = 30
$ws.Cells.Item(133, 11).Value = 15000
$ws.Cells.Item(133, 12).Value = 15000
$ws.Cells.Item(133, 13).Value = 15000
$ws.Cells.Item(133, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(133, 15).Value = 'Región del Maule'
$ws.Cells.Item(133, 16).Value = 600
$ws.Cells.Item(133, 17).Value = 25

# Row 134
$ws.Cells.Item(134, 4).Value = 44643
$ws.Cells.Item(134, 8).Value = 'Americana (o)'
$ws.Cells.Item(134, 10).Value = 60
$ws.Cells.Item(134, 11).Value = 8500
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = 8750
$ws.Cells.Item(134, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(134, 15).Value = 'Región del Maule'
$ws.Cells.Item(134, 16).Value = 583
$ws.Cells.Item(134, 17).Value = 15

# Row 135
$ws.Cells.Item(135, 4).Value = 45005
$ws.Cells.Item(135, 8).Value = 'Cacho cabra rojo'
$ws.Cells.Item(135, 10).Value = 40
$ws.Cells.Item(135, 11).Value = 15000
$ws.Cells.Item(135, 12).Value = 15000
$ws.Cells.Item(135, 13).Value = 15000
$ws.Cells.Item(135, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(135, 15).Value = 'Región del Maule'
$ws.Cells.Item(135, 16).Value = 600
$ws.Cells.Item(135, 17).Value = 25

# Row 136
$ws.Cells.Item(136, 4).Value = 45005
$ws.Cells.Item(136, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(136, 10).Value = 30
$ws.Cells.Item(136, 11).Value = 15000
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = 15000
$ws.Cells.Item(136, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(136, 15).Value = 'Región del Maule'
$ws.Cells.Item(136, 16).Value = 600
$ws.Cells.Item(136, 17).Value = 25

# Row 137
$ws.Cells.Item(137, 4).Value = 45005
$ws.Cells.Item(137, 8).Value = 'Cristal'
$ws.Cells.Item(137, 10).Value = 20
$ws.Cells.Item(137, 11).Value = 15000
$ws.Cells.Item(137, 12).Value = 15000
$ws.Cells.Item(137, 13).Value = 15000
$ws.Cells.Item(137, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(137, 15).Value = 'Región del Maule'
$ws.Cells.Item(137, 16).Value = 600
$ws.Cells.Item(137, 17).Value = 25

# Row 138
$ws.Cells.Item(138, 4).Value = 44901
$ws.Cells.Item(138, 8).Value = 'Americana (o)'
$ws.Cells.Item(138, 10).Value = 60
$ws.Cells.Item(138, 11).Value = 16000
$ws.Cells.Item(138, 12).Value = 17000
$ws.Cells.Item(138, 13).Value = 16500
$ws.Cells.Item(138, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(138, 15).Value = 'Región del Maule'
$ws.Cells.Item(138, 16).Value = 1100
$ws.Cells.Item(138, 17).Value = 15

# Row 139
$ws.Cells.Item(139, 4).Value = 44859
$ws.Cells.Item(139, 8).Value = 'Inferno'
$ws.Cells.Item(139, 10).Value = 50
$ws.Cells.Item(139, 11).Value = 20000
$ws.Cells.Item(139, 12).Value = 20000
$ws.Cells.Item(139, 13).Value = 20000
$ws.Cells.Item(139, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(139, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(139, 16).Value = 2000
$ws.Cells.Item(139, 17).Value = 10

# Row 140
$ws.Cells.Item(140, 4).Value = 44873
$ws.Cells.Item(140, 8).Value = 'Americana (o)'
$ws.Cells.Item(140, 10).Value = 50
$ws.Cells.Item(140, 11).Value = 35000
$ws.Cells.Item(140, 12).Value = 35000
$ws.Cells.Item(140, 13).Value = 35000
$ws.Cells.Item(140, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(140, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(140, 16).Value = 1400
$ws.Cells.Item(140, 17).Value = 25

# Row 141
$ws.Cells.Item(141, 4).Value = 44873
$ws.Cells.Item(141, 8).Value = 'Inferno'
$ws.Cells.Item(141, 10).Value = 60
$ws.Cells.Item(141, 11).Value = 20000
$ws.Cells.Item(141, 12).Value = 22000
$ws.Cells.Item(141, 13).Value = 21000
$ws.Cells.Item(141, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(141, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(141, 16).Value = 2100
$ws.Cells.Item(141, 17).Value = 10

# Row 142
$ws.Cells.Item(142, 4).Value = 44238
$ws.Cells.Item(142, 8).Value = 'Americana (o)'
$ws.Cells.Item(142, 10).Value = 30
$ws.Cells.Item(142, 11).Value = 14000
$ws.Cells.Item(142, 12).Value = 15000
$ws.Cells.Item(142, 13).Value = 14500
$ws.Cells.Item(142, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(142, 15).Value = 'Región del Maule'
$ws.Cells.Item(142, 16).Value = 580
$ws.Cells.Item(142, 17).Value = 25

# Row 143
$ws.Cells.Item(143, 4).Value = 44165
$ws.Cells.Item(143, 8).Value = 'Americana (o)'
$ws.Cells.Item(143, 10).Value = 75
$ws.Cells.Item(143, 11).Value = 2000
$ws.Cells.Item(143, 12).Value = 2000
$ws.Cells.Item(143, 13).Value = 2000
$ws.Cells.Item(143, 14).Value = '$/kilo'
$ws.Cells.Item(143, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(143, 16).Value = 2000
$ws.Cells.Item(143, 17).Value = 1

# Row 144
$ws.Cells.Item(144, 4).Value = 44204
$ws.Cells.Item(144, 8).Value = 'Americana (o)'
$ws.Cells.Item(144, 10).Value = 75
$ws.Cells.Item(144, 11).Value = 1200
$ws.Cells.Item(144, 12).Value = 1400
$ws.Cells.Item(144, 13).Value = 1267
$ws.Cells.Item(144, 14).Value = '$/kilo'
$ws.Cells.Item(144, 15).Value = 'Región del Maule'
$ws.Cells.Item(144, 16).Value = 1267
$ws.Cells.Item(144, 17).Value = 1

# Row 145
$ws.Cells.Item(145, 4).Value = 44229
$ws.Cells.Item(145, 8).Value = 'Americana (o)'
$ws.Cells.Item(145, 10).Value = 30
$ws.Cells.Item(145, 11).Value = 12000
$ws.Cells.Item(145, 12).Value = 13000
$ws.Cells.Item(145, 13).Value = 12500
$ws.Cells.Item(145, 14).Value = '$/caja 14 kilos'
$ws.Cells.Item(145, 15).Value = 'Región del Maule'
$ws.Cells.Item(145, 16).Value = 893
$ws.Cells.Item(145, 17).Value = 14

# Row 146
$ws.Cells.Item(146, 4).Value = 44804
$ws.Cells.Item(146, 8).Value = 'Inferno'
$ws.Cells.Item(146, 10).Value = 60
$ws.Cells.Item(146, 11).Value = 17000
$ws.Cells.Item(146, 12).Value = 18000
$ws.Cells.Item(146, 13).Value = 17500
$ws.Cells.Item(146, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(146, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(146, 16).Value = 1167
$ws.Cells.Item(146, 17).Value = 15

# Row 147
$ws.Cells.Item(147, 4).Value = 44663
$ws.Cells.Item(147, 8).Value = 'Americana (o)'
$ws.Cells.Item(147, 10).Value = 80
$ws.Cells.Item(147, 11).Value = 8500
$ws.Cells.Item(147, 12).Value = 9000
$ws.Cells.Item(147, 13).Value = 8750
$ws.Cells.Item(147, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(147, 15).Value = 'Región del Maule'
$ws.Cells.Item(147, 16).Value = 583
$ws.Cells.Item(147, 17).Value = 15

# Row 148
$ws.Cells.Item(148, 4).Value = 44285
$ws.Cells.Item(148, 8).Value = 'Cristal'
$ws.Cells.Item(148, 10).Value = 40
$ws.Cells.Item(148, 11).Value = 14000
$ws.Cells.Item(148, 12).Value = 15000
$ws.Cells.Item(148, 13).Value = 14500
$ws.Cells.Item(148, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(148, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(148, 16).Value = 580
$ws.Cells.Item(148, 17).Value = 25

# Row 149
$ws.Cells.Item(149, 4).Value = 44292
$ws.Cells.Item(149, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(149, 10).Value = 60
$ws.Cells.Item(149, 11).Value = 13000
$ws.Cells.Item(149, 12).Value = 14000
$ws.Cells.Item(149, 13).Value = 13500
$ws.Cells.Item(149, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(149, 15).Value = 'Región del Maule'
$ws.Cells.Item(149, 16).Value = 540
$ws.Cells.Item(149, 17).Value = 25

# Row 150
$ws.Cells.Item(150, 4).Value = 44910
$ws.Cells.Item(150, 8).Value = 'Americana (o)'
$ws.Cells.Item(150, 10).Value = 50
$ws.Cells.Item(150, 11).Value = 15000
$ws.Cells.Item(150, 12).Value = 15000
$ws.Cells.Item(150, 13).Value = 15000
$ws.Cells.Item(150, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(150, 15).Value = 'Región del Maule'
$ws.Cells.Item(150, 16).Value = 1000
$ws.Cells.Item(150, 17).Value = 15

# Row 151
$ws.Cells.Item(151, 4).Value = 44202
$ws.Cells.Item(151, 8).Value = 'Americana (o)'
$ws.Cells.Item(151, 10).Value = 75
$ws.Cells.Item(151, 11).Value = 1500
$ws.Cells.Item(151, 12).Value = 1700
$ws.Cells.Item(151, 13).Value = 1567
$ws.Cells.Item(151, 14).Value = '$/kilo'
$ws.Cells.Item(151, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(151, 16).Value = 1567
$ws.Cells.Item(151, 17).Value = 1

# Row 152
$ws.Cells.Item(152, 4).Value = 44651
$ws.Cells.Item(152, 8).Value = 'Americana (o)'
$ws.Cells.Item(152, 10).Value = 60
$ws.Cells.Item(152, 11).Value = 8500
$ws.Cells.Item(152, 12).Value = 9000
$ws.Cells.Item(152, 13).Value = 8750
$ws.Cells.Item(152, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(152, 15).Value = 'Región del Maule'
$ws.Cells.Item(152, 16).Value = 583
$ws.Cells.Item(152, 17).Value = 15

# Row 153
$ws.Cells.Item(153, 4).Value = 44571
$ws.Cells.Item(153, 8).Value = 'Americana (o)'
$ws.Cells.Item(153, 10).Value = 60
$ws.Cells.Item(153, 11).Value = 15000
$ws.Cells.Item(153, 12).Value = 16000
$ws.Cells.Item(153, 13).Value = 15500
$ws.Cells.Item(153, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(153, 15).Value = 'Región del Maule'
$ws.Cells.Item(153, 16).Value = 1033
$ws.Cells.Item(153, 17).Value = 15

# Row 154
$ws.Cells.Item(154, 4).Value = 44860
$ws.Cells.Item(154, 8).Value = 'Americana (o)'
$ws.Cells.Item(154, 10).Value = 30
$ws.Cells.Item(154, 11).Value = 45000
$ws.Cells.Item(154, 12).Value = 45000
$ws.Cells.Item(154, 13).Value = 45000
$ws.Cells.Item(154, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(154, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(154, 16).Value = 1800
$ws.Cells.Item(154, 17).Value = 25

# Row 155
$ws.Cells.Item(155, 4).Value = 44939
$ws.Cells.Item(155, 8).Value = 'Americana (o)'
$ws.Cells.Item(155, 10).Value = 30
$ws.Cells.Item(155, 11).Value = 13000
$ws.Cells.Item(155, 12).Value = 13000
$ws.Cells.Item(155, 13).Value = 13000
$ws.Cells.Item(155, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(155, 15).Value = 'Región del Maule'
$ws.Cells.Item(155, 16).Value = 867
$ws.Cells.Item(155, 17).Value = 15

# Row 156
$ws.Cells.Item(156, 4).Value = 44673
$ws.Cells.Item(156, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(156, 10).Value = 30
$ws.Cells.Item(156, 11).Value = 16000
$ws.Cells.Item(156, 12).Value = 17000
$ws.Cells.Item(156, 13).Value = 16500
$ws.Cells.Item(156, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(156, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(156, 16).Value = 660
$ws.Cells.Item(156, 17).Value = 25

# Row 157
$ws.Cells.Item(157, 4).Value = 44568
$ws.Cells.Item(157, 8).Value = 'Americana (o)'
$ws.Cells.Item(157, 10).Value = 100
$ws.Cells.Item(157, 11).Value = 19000
$ws.Cells.Item(157, 12).Value = 20000
$ws.Cells.Item(157, 13).Value = 19500
$ws.Cells.Item(157, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(157, 15).Value = 'Región del Maule'
$ws.Cells.Item(157, 16).Value = 1300
$ws.Cells.Item(157, 17).Value = 15

# Row 158
$ws.Cells.Item(158, 4).Value = 44874
$ws.Cells.Item(158, 8).Value = 'Americana (o)'
$ws.Cells.Item(158, 10).Value = 30
$ws.Cells.Item(158, 11).Value = 35000
$ws.Cells.Item(158, 12).Value = 35000
$ws.Cells.Item(158, 13).Value = 35000
$ws.Cells.Item(158, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(158, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(158, 16).Value = 1400
$ws.Cells.Item(158, 17).Value = 25

# Row 159
$ws.Cells.Item(159, 4).Value = 44620
$ws.Cells.Item(159, 8).Value = 'Americana (o)'
$ws.Cells.Item(159, 10).Value = 30
$ws.Cells.Item(159, 11).Value = 9000
$ws.Cells.Item(159, 12).Value = 9000
$ws.Cells.Item(159, 13).Value = 9000
$ws.Cells.Item(159, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(159, 15).Value = 'Región del Maule'
$ws.Cells.Item(159, 16).Value = 600
$ws.Cells.Item(159, 17).Value = 15

# Row 160
$ws.Cells.Item(160, 4).Value = 44257
$ws.Cells.Item(160, 8).Value = 'Americana (o)'
$ws.Cells.Item(160, 10).Value = 27
$ws.Cells.Item(160, 11).Value = 15000
$ws.Cells.Item(160, 12).Value = 16000
$ws.Cells.Item(160, 13).Value = 15556
$ws.Cells.Item(160, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(160, 15).Value = 'Región del Maule'
$ws.Cells.Item(160, 16).Value = 622
$ws.Cells.Item(160, 17).Value = 25

# Row 161
$ws.Cells.Item(161, 4).Value = 45008
$ws.Cells.Item(161, 8).Value = 'Cacho cabra rojo'
$ws.Cells.Item(161, 10).Value = 25
$ws.Cells.Item(161, 11).Value = 15000
$ws.Cells.Item(161, 12).Value = 15000
$ws.Cells.Item(161, 13).Value = 15000
$ws.Cells.Item(161, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(161, 15).Value = 'Región del Maule'
$ws.Cells.Item(161, 16).Value = 600
$ws.Cells.Item(161, 17).Value = 25

# Row 162
$ws.Cells.Item(162, 4).Value = 45008
$ws.Cells.Item(162, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(162, 10).Value = 20
$ws.Cells.Item(162, 11).Value = 15000
$ws.Cells.Item(162, 12).Value = 15000
$ws.Cells.Item(162, 13).Value = 15000
$ws.Cells.Item(162, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(162, 15).Value = 'Región del Maule'
$ws.Cells.Item(162, 16).Value = 600
$ws.Cells.Item(162, 17).Value = 25

# Row 163
$ws.Cells.Item(163, 4).Value = 45008
$ws.Cells.Item(163, 8).Value = 'Cristal'
$ws.Cells.Item(163, 10).Value = 20
$ws.Cells.Item(163, 11).Value = 15000
$ws.Cells.Item(163, 12).Value = 15000
$ws.Cells.Item(163, 13).Value = 15000
$ws.Cells.Item(163, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(163, 15).Value = 'Región del Maule'
$ws.Cells.Item(163, 16).Value = 600
$ws.Cells.Item(163, 17).Value = 25

# Row 164
$ws.Cells.Item(164, 4).Value = 44960
$ws.Cells.Item(164, 8).Value = 'Americana (o)'
$ws.Cells.Item(164, 10).Value = 60
$ws.Cells.Item(164, 11).Value = 11000
$ws.Cells.Item(164, 12).Value = 12000
$ws.Cells.Item(164, 13).Value = 11500
$ws.Cells.Item(164, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(164, 15).Value = 'Región del Maule'
$ws.Cells.Item(164, 16).Value = 767
$ws.Cells.Item(164, 17).Value = 15

# Row 165
$ws.Cells.Item(165, 4).Value = 44771
$ws.Cells.Item(165, 8).Value = 'Americana (o)'
$ws.Cells.Item(165, 10).Value = 30
$ws.Cells.Item(165, 11).Value = 35000
$ws.Cells.Item(165, 12).Value = 35000
$ws.Cells.Item(165, 13).Value = 35000
$ws.Cells.Item(165, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(165, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(165, 16).Value = 1400
$ws.Cells.Item(165, 17).Value = 25

# Row 166
$ws.Cells.Item(166, 4).Value = 44972
$ws.Cells.Item(166, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(166, 10).Value = 30
$ws.Cells.Item(166, 11).Value = 15000
$ws.Cells.Item(166, 12).Value = 15000
$ws.Cells.Item(166, 13).Value = 15000
$ws.Cells.Item(166, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(166, 15).Value = 'Región del Maule'
$ws.Cells.Item(166, 16).Value = 600
$ws.Cells.Item(166, 17).Value = 25

# Row 167
$ws.Cells.Item(167, 4).Value = 44988
$ws.Cells.Item(167, 8).Value = 'Cacho cabra verde'
$ws.Cells.Item(167, 10).Value = 30
$ws.Cells.Item(167, 11).Value = 14000
$ws.Cells.Item(167, 12).Value = 14000
$ws.Cells.Item(167, 13).Value = 14000
$ws.Cells.Item(167, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(167, 15).Value = 'Región del Maule'
$ws.Cells.Item(167, 16).Value = 560
$ws.Cells.Item(167, 17).Value = 25

# Row 168
$ws.Cells.Item(168, 4).Value = 44608
$ws.Cells.Item(168, 8).Value = 'Americana (o)'
$ws.Cells.Item(168, 10).Value = 100
$ws.Cells.Item(168, 11).Value = 9000
$ws.Cells.Item(168, 12).Value = 9500
$ws.Cells.Item(168, 13).Value = 9250
$ws.Cells.Item(168, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(168, 15).Value = 'Región del Maule'
$ws.Cells.Item(168, 16).Value = 617
$ws.Cells.Item(168, 17).Value = 15

# Row 169
$ws.Cells.Item(169, 4).Value = 44925
$ws.Cells.Item(169, 8).Value = 'Americana (o)'
$ws.Cells.Item(169, 10).Value = 60
$ws.Cells.Item(169, 11).Value = 13000
$ws.Cells.Item(169, 12).Value = 14000
$ws.Cells.Item(169, 13).Value = 13500
$ws.Cells.Item(169, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(169, 15).Value = 'Región del Maule'
$ws.Cells.Item(169, 16).Value = 900
$ws.Cells.Item(169, 17).Value = 15

# Make sure the date column keeps its date format for the newly added rows too
$ws.Range("D167:D169").NumberFormat = $ws.Range("D166").NumberFormat

Write-Output "Done. Used range now covers through row $($ws.Cells.Item(169,18).Row)"
